$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update masthead volume/number and week-covering date range text (merged header cells)
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# Column E narrows to match the other % Chg columns
$ws.Columns("E").ColumnWidth = 6.168446

# --- Cells switching from a numeric format to text (blank-count / undefined-% placeholders) ---
# Pull number-format/style from a stable style-14 (text) cell, e.g. G15, onto the target,
# then write the placeholder text (leading apostrophe forces text even though it looks numeric).
$ws.Range("F15").Value = "'0"
$ws.Range("G15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("D16").Value = "'0"
$ws.Range("G15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "'***.*"
$ws.Range("G15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D18").Value = "'0"
$ws.Range("G15").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "'***.*"
$ws.Range("G15").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("C31").Value = "'0"
$ws.Range("G15").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D31").Value = "'0"
$ws.Range("G15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$ws.Range("G15").Copy()
$ws.Range("E31").PasteSpecial(-4122)

# --- Cells switching from text placeholders back to numeric counts ---
$ws.Range("D20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1
$ws.Range("D20").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 61
$ws.Range("K16").Value = 60.526315789473
$ws.Range("L16").Value = 15.094339622641
$ws.Range("M16").Value = 69.444444444444
$ws.Range("N16").Value = -83.989501312336
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 225
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 65
$ws.Range("I17").Value = 76
$ws.Range("J17").Value = 61
$ws.Range("K17").Value = 24.590163934426
$ws.Range("L17").Value = 49.019607843137
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = -48.648648648648
$ws.Range("C18").Value = 7
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 56
$ws.Range("K18").Value = 33.333333333333
$ws.Range("L18").Value = -23.287671232876
$ws.Range("M18").Value = -6.666666666666
$ws.Range("N18").Value = -91.437308868501
$ws.Range("C19").Value = 35
$ws.Range("D19").Value = 41
$ws.Range("E19").Value = -14.634146341463
$ws.Range("F19").Value = 134
$ws.Range("H19").Value = -4.285714285714
$ws.Range("I19").Value = 534
$ws.Range("J19").Value = 552
$ws.Range("K19").Value = -3.260869565217
$ws.Range("L19").Value = 4.093567251461
$ws.Range("M19").Value = 15.83514099783
$ws.Range("N19").Value = -74.595623215984
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 12
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -47.826086956521
$ws.Range("L20").Value = -58.620689655172
$ws.Range("M20").Value = 9.090909090909
$ws.Range("N20").Value = -92.207792207792
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 20.833333333333
$ws.Range("F21").Value = 195
$ws.Range("G21").Value = 190
$ws.Range("H21").Value = 2.631578947368
$ws.Range("I21").Value = 743
$ws.Range("J21").Value = 718
$ws.Range("K21").Value = 3.481894150417
$ws.Range("L21").Value = 2.482758620689
$ws.Range("M21").Value = 20.421393841166
$ws.Range("N21").Value = -78.501157407407
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 18.181818181818
$ws.Range("L22").Value = 30
$ws.Range("M22").Value = 36.842105263157
$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = 8.163265306122
$ws.Range("F24").Value = 180
$ws.Range("G24").Value = 194
$ws.Range("H24").Value = -7.21649484536
$ws.Range("I24").Value = 849
$ws.Range("J24").Value = 743
$ws.Range("K24").Value = 14.266487213997
$ws.Range("L24").Value = 34.976152623211
$ws.Range("M24").Value = 64.53488372093
$ws.Range("C25").Value = 61
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = 24.489795918367
$ws.Range("F25").Value = 183
$ws.Range("G25").Value = 185
$ws.Range("H25").Value = -1.081081081081
$ws.Range("I25").Value = 848
$ws.Range("J25").Value = 773
$ws.Range("K25").Value = 9.702457956015
$ws.Range("L25").Value = 23.435225618631
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 110
$ws.Range("F26").Value = 56
$ws.Range("H26").Value = -5.084745762711
$ws.Range("I26").Value = 197
$ws.Range("J26").Value = 195
$ws.Range("K26").Value = 1.025641025641
$ws.Range("L26").Value = 37.762237762237
$ws.Range("M26").Value = 33.108108108108
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -45.454545454545
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 10
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 22
$ws.Range("J28").Value = 22
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -24.137931034482
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = -75
$ws.Range("L31").Value = 33.333333333333

Write-Output "done"